# Update test manifest data for sample ID change: CMO_SAMPLE_ID values gain an
# "_IGO" suffix (fixing up the stray "test_sample_1_N_IGO_TEST" id along the
# way), INVESTIGATOR_SAMPLE_ID values are populated, and the SampleRenames
# sheet is updated to map each new IGO id back to its clean base name.

$wb = $excel.ActiveWorkbook
$wsInfo = $wb.Worksheets.Item("SampleInfo")
$wsRenames = $wb.Worksheets.Item("SampleRenames")

# New CMO_SAMPLE_ID (column A), INVESTIGATOR_SAMPLE_ID (column C) and the
# clean base name (used on the SampleRenames sheet) for each data row.
$rowData = @{
    2 = @{ CmoId = "test_sample_2_T_IGO";    Investigator = "test_investigator_sample_2_T"; Base = "test_sample_2_T" }
    3 = @{ CmoId = "test_sample_1_N_IGO";    Investigator = "test_investigator_sample_1_N"; Base = "test_sample_1_N" }
    4 = @{ CmoId = "test_sample_4_T_IGO";    Investigator = "test_investigator_sample_4_T"; Base = "test_sample_4_T" }
    5 = @{ CmoId = "test_sample_3_N_IGO";    Investigator = "test_investigator_sample_3_N"; Base = "test_sample_3_N" }
    6 = @{ CmoId = "test_sample_6_T_IGO";    Investigator = "test_investigator_sample_6_T"; Base = "test_sample_6_T" }
    7 = @{ CmoId = "test_sample_5_N_IGO";    Investigator = "test_investigator_sample_5_N"; Base = "test_sample_5_N" }
}

for ($row = 2; $row -le 7; $row++) {
    $data = $rowData[$row]

    $wsInfo.Cells.Item($row, 1).Value = $data.CmoId
    $wsInfo.Cells.Item($row, 3).Value = $data.Investigator

    # Update SampleRenames sheet accordingly: column A becomes the new
    # (IGO-suffixed) name and column B becomes the original base name.
    $wsRenames.Cells.Item($row, 1).Value = $data.CmoId
    $wsRenames.Cells.Item($row, 2).Value = $data.Base
}

# Column A on SampleRenames no longer carries the special style (s="1"),
# revert it to the default style while column B keeps its existing style.
$wsRenames.Range("A2:A7").Style = "Normal"

# Update the recorded selection (active cell) on each sheet, leaving
# SampleRenames as the active sheet/tab, matching the saved workbook state.
[void]$wsInfo.Activate()
[void]$wsInfo.Range("I11").Select()
[void]$wsRenames.Activate()
[void]$wsRenames.Range("C25").Select()
